$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

# Reset cell style to default ("Normal") before writing the new value so
# the explicit (green-fill) style is dropped, then update the mobile
# number value.
$ws.Range("B5").Style = "Normal"
$ws.Range("B5").Value = 9912345678

# Make "Input" the active sheet/tab and move the selection to B5.
$ws.Activate()
$ws.Range("B5").Select()
